# Natmi following Dr Hou advice
# Update the Wnt2-Fzd9 LR-pair sheet with recomputed detection/expression
# and specificity values for all six "Target cluster" rows (rows 2-7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns:
#  E = Ligand-expressing cells
#  G = Ligand average expression value
#  H = Ligand total expression value
#  K = Receptor-expressing cells
#  M = Receptor average expression value
#  N = Receptor total expression value
#  O = Receptor derived specificity of average expression value
#  P = Receptor derived specificity of total expression value
#  Q = Edge average expression weight
#  R = Edge total expression weight
#  S = Edge average expression derived specificity
#  T = Edge total expression derived specificity

$rows = @{
    2 = @{ E = 3; G = 0.574538; H = 1.723614; K = 2; M = 0.7517265; N = 1.503453;
           O = 0.3430414969595709; P = 0.2862361877440628; Q = 0.431895439857; R = 2.591372639142;
           S = 0.3430414969595709; T = 0.2862361877440628 }
    3 = @{ E = 3; G = 0.574538; H = 1.723614; K = 3; M = 0.2592623333333333; N = 0.777787;
           O = 0.1183112992982127; P = 0.1480796444962971; Q = 0.1489560624686667; R = 1.340604562218;
           S = 0.1183112992982127; T = 0.1480796444962971 }
    4 = @{ E = 3; G = 0.574538; H = 1.723614; K = 3; M = 0.1423686666666667; N = 0.427106;
           O = 0.06496825711674591; P = 0.08131494180570706; Q = 0.08179620900933333; R = 0.7361658810839999;
           S = 0.06496825711674591; T = 0.08131494180570706 }
    5 = @{ E = 3; G = 0.574538; H = 1.723614; K = 3; M = 0.171105; N = 0.513315;
           O = 0.07808174294409917; P = 0.09772791614493011; Q = 0.09830632448999999; R = 0.88475692041;
           S = 0.07808174294409917; T = 0.09772791614493011 }
    6 = @{ E = 3; G = 0.574538; H = 1.723614; K = 3; M = 0.2970403333333334; N = 0.8911210000000001;
           O = 0.1355508427653363; P = 0.1696568352044773; Q = 0.1706609590326667; R = 1.535948631294;
           S = 0.1355508427653363; T = 0.1696568352044773 }
    7 = @{ E = 3; G = 0.574538; H = 1.723614; K = 2; M = 0.5698545; N = 1.139709;
           O = 0.260046360916035; P = 0.2169844746045258; Q = 0.327403064721; R = 1.964418388326;
           S = 0.260046360916035; T = 0.2169844746045258 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("E$r").Value = $vals.E
    $ws.Range("G$r").Value = $vals.G
    $ws.Range("H$r").Value = $vals.H
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("Q$r").Value = $vals.Q
    $ws.Range("R$r").Value = $vals.R
    $ws.Range("S$r").Value = $vals.S
    $ws.Range("T$r").Value = $vals.T
}
